$d = $word.ActiveDocument

# Update the title/date paragraph
$d.Content.Find.Execute("2025-03-27 Thursday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-03-28 Friday", 2) | Out-Null

# Update each table cell value positionally (row-major order)
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "94-66="
$t.Cell(1, 2).Range.Text = "85-77="
$t.Cell(1, 3).Range.Text = "19+53="
$t.Cell(1, 4).Range.Text = "24-15="
$t.Cell(1, 5).Range.Text = "21-17="
$t.Cell(2, 1).Range.Text = "45+28="
$t.Cell(2, 2).Range.Text = "18+68="
$t.Cell(2, 3).Range.Text = "84-17="
$t.Cell(2, 4).Range.Text = "92-55="
$t.Cell(2, 5).Range.Text = "18+14="
$t.Cell(3, 1).Range.Text = "64-39="
$t.Cell(3, 2).Range.Text = "47+29="
$t.Cell(3, 3).Range.Text = "14-6="
$t.Cell(3, 4).Range.Text = "94-88="
$t.Cell(3, 5).Range.Text = "81-49="
$t.Cell(4, 1).Range.Text = "57-29="
$t.Cell(4, 2).Range.Text = "22+9="
$t.Cell(4, 3).Range.Text = "15+7="
$t.Cell(4, 4).Range.Text = "96-8="
$t.Cell(4, 5).Range.Text = "96-8="
$t.Cell(5, 1).Range.Text = "27+9="
$t.Cell(5, 2).Range.Text = "79+17="
$t.Cell(5, 3).Range.Text = "94-89="
$t.Cell(5, 4).Range.Text = "6+37="
$t.Cell(5, 5).Range.Text = "8+77="
$t.Cell(6, 1).Range.Text = "16+75="
$t.Cell(6, 2).Range.Text = "70-66="
$t.Cell(6, 3).Range.Text = "27+58="
$t.Cell(6, 4).Range.Text = "73+9="
$t.Cell(6, 5).Range.Text = "95-58="
$t.Cell(7, 1).Range.Text = "29+19="
$t.Cell(7, 2).Range.Text = "80-13="
$t.Cell(7, 3).Range.Text = "82-45="
$t.Cell(7, 4).Range.Text = "35+57="
$t.Cell(7, 5).Range.Text = "72-69="
$t.Cell(8, 1).Range.Text = "7+38="
$t.Cell(8, 2).Range.Text = "92-65="
$t.Cell(8, 3).Range.Text = "66+5="
$t.Cell(8, 4).Range.Text = "28+4="
$t.Cell(8, 5).Range.Text = "95-38="
$t.Cell(9, 1).Range.Text = "56-28="
$t.Cell(9, 2).Range.Text = "48+49="
$t.Cell(9, 3).Range.Text = "95-86="
$t.Cell(9, 4).Range.Text = "29+47="
$t.Cell(9, 5).Range.Text = "28+44="
$t.Cell(10, 1).Range.Text = "71-4="
$t.Cell(10, 2).Range.Text = "79+6="
$t.Cell(10, 3).Range.Text = "45+37="
$t.Cell(10, 4).Range.Text = "82-46="
$t.Cell(10, 5).Range.Text = "63-46="
$t.Cell(11, 1).Range.Text = "49+8="
$t.Cell(11, 2).Range.Text = "60-53="
$t.Cell(11, 3).Range.Text = "53-29="
$t.Cell(11, 4).Range.Text = "91-24="
$t.Cell(11, 5).Range.Text = "63-36="
$t.Cell(12, 1).Range.Text = "29+46="
$t.Cell(12, 2).Range.Text = "36+27="
$t.Cell(12, 3).Range.Text = "55-38="
$t.Cell(12, 4).Range.Text = "68-59="
$t.Cell(12, 5).Range.Text = "67+6="
$t.Cell(13, 1).Range.Text = "70-61="
$t.Cell(13, 2).Range.Text = "24+9="
$t.Cell(13, 3).Range.Text = "7+26="
$t.Cell(13, 4).Range.Text = "19+2="
$t.Cell(13, 5).Range.Text = "19+37="
$t.Cell(14, 1).Range.Text = "33+49="
$t.Cell(14, 2).Range.Text = "56+27="
$t.Cell(14, 3).Range.Text = "92-18="
$t.Cell(14, 4).Range.Text = "14+77="
$t.Cell(14, 5).Range.Text = "16+29="
$t.Cell(15, 1).Range.Text = "59+38="
$t.Cell(15, 2).Range.Text = "7+25="
$t.Cell(15, 3).Range.Text = "90-52="
$t.Cell(15, 4).Range.Text = "94-18="
$t.Cell(15, 5).Range.Text = "59+35="
$t.Cell(16, 1).Range.Text = "48+13="
$t.Cell(16, 2).Range.Text = "85-36="
$t.Cell(16, 3).Range.Text = "27+17="
$t.Cell(16, 4).Range.Text = "72-24="
$t.Cell(16, 5).Range.Text = "63+18="
$t.Cell(17, 1).Range.Text = "64+28="
$t.Cell(17, 2).Range.Text = "82-77="
$t.Cell(17, 3).Range.Text = "80-21="
$t.Cell(17, 4).Range.Text = "19+44="
$t.Cell(17, 5).Range.Text = "56+39="
$t.Cell(18, 1).Range.Text = "7+47="
$t.Cell(18, 2).Range.Text = "72-38="
$t.Cell(18, 3).Range.Text = "71-38="
$t.Cell(18, 4).Range.Text = "32+59="
$t.Cell(18, 5).Range.Text = "60-23="
$t.Cell(19, 1).Range.Text = "90-46="
$t.Cell(19, 2).Range.Text = "55+18="
$t.Cell(19, 3).Range.Text = "28+25="
$t.Cell(19, 4).Range.Text = "16+77="
$t.Cell(19, 5).Range.Text = "22-13="
$t.Cell(20, 1).Range.Text = "96-9="
$t.Cell(20, 2).Range.Text = "93-88="
$t.Cell(20, 3).Range.Text = "26+18="
$t.Cell(20, 4).Range.Text = "16+65="
$t.Cell(20, 5).Range.Text = "31-6="
